# Split the run containing "{m" (the opening brace and the 'm' of the {m:...}
# field marker) into two separate runs: "{" and "m". This mirrors the switch
# to TokenIteratorFieldRewriterSplit, which now emits the field-open
# character and the following token text as distinct runs.
#
# Likewise split the run containing ".setWidth(100)}" into two separate
# runs: ".setWidth(100)" (keeping its original orange run formatting) and
# "}" (the closing field marker, emitted as its own run).

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(2)
$pStart = $p.Range.Start

# ---------------------------------------------------------------------
# 1) "{m" -> "{" + "m"
# ---------------------------------------------------------------------
$paraText = $p.Range.Text
$openIdx = $paraText.IndexOf("{m")
$openStart = $pStart + $openIdx
$openEnd = $openStart + 2

$openRun = $d.Range($openStart, $openEnd)

# Delete the combined "{m" text, then insert "{" and "m" as two distinct
# runs at the now-collapsed insertion point so they are not recombined
# into a single run.
$openRun.Text = ""
$openInsertionPoint = $d.Range($openStart, $openStart)
$openXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$openInsertionPoint.InsertXML($openXml)

# ---------------------------------------------------------------------
# 2) ".setWidth(100)}" -> ".setWidth(100)" + "}"
# ---------------------------------------------------------------------
$paraText2 = $p.Range.Text
$closeIdx = $paraText2.IndexOf(".setWidth(100)}")
$closeStart = $pStart + $closeIdx
$closeEnd = $closeStart + ".setWidth(100)}".Length

$closeRun = $d.Range($closeStart, $closeEnd)
$closeXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>.setWidth(100)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$closeRun.InsertXML($closeXml)
